$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MPWR")

# Row 4 - Inventory
$ws.Range("B4").Value = 157000000.0
$ws.Range("C4").Value = 148000000.0
$ws.Range("D4").Value = 152000000.0
$ws.Range("E4").Value = 131000000.0
$ws.Range("F4").Value = 128000000.0

# Row 14 - Accounts Payable
$ws.Range("B14").Value = 38000000.0
$ws.Range("C14").Value = 48000000.0
$ws.Range("D14").Value = 45000000.0
$ws.Range("E14").Value = 38000000.0
$ws.Range("F14").Value = 27000000.0

# Row 20 - Long Term Tax Liability (Deferred)
$ws.Range("B20").Value = -19000000.0
$ws.Range("D20").Value = -13000000.0
$ws.Range("E20").Value = -14000000.0
$ws.Range("F20").Value = -17000000.0
